$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original (date-only) number format used by row 23 before
# changing it, so the new row 24 can reuse it.
$dateOnlyFormat = $ws.Range("A23").NumberFormat

# Row 23, column A switches to the datetime format used by the rest of the
# column (e.g. A22), matching the edit where only the newest row keeps the
# date-only formatting.
$ws.Range("A23").NumberFormat = $ws.Range("A22").NumberFormat

# New row 24 with the next day's raw + clean SSA data (June 23rd).
$ws.Range("A24").Value = 44005
$ws.Range("B24").Value = 191410
$ws.Range("C24").Value = 251355
$ws.Range("D24").Value = 59106
$ws.Range("E24").Value = 23377
$ws.Range("F24").Value = 31.6

# A24 takes on the date-only format row 23 used to have.
$ws.Range("A24").NumberFormat = $dateOnlyFormat
